$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 194.4
$ws.Range("I8").Value = 194.4
$ws.Range("K8").Value = 583.2
$ws.Range("M8").Value = -444.2
$ws.Range("H17").Value = 471.33334
$ws.Range("J17").Value = 471.33334
$ws.Range("L17").Value = 1414.00002
$ws.Range("N17").Value = -1750.00002
$ws.Range("H132").Value = 1553.5
$ws.Range("I132").Value = 1553.5
$ws.Range("K132").Value = 4660.5
$ws.Range("M132").Value = -2130.5
$ws.Range("H135").Value = 3609.7144
$ws.Range("I135").Value = 4566.4
$ws.Range("J135").Value = 1218
$ws.Range("K135").Value = 41097.6
$ws.Range("L135").Value = 10962
$ws.Range("M135").Value = -38562.6
$ws.Range("N135").Value = -16032
$ws.Range("H137").Value = 7137.3335
$ws.Range("I137").Value = 30449
$ws.Range("J137").Value = 4683.4736
$ws.Range("K137").Value = 91347
$ws.Range("L137").Value = 14050.4208
$ws.Range("M137").Value = -88797
$ws.Range("N137").Value = -19150.4208
$ws.Range("H138").Value = 3820.9524
$ws.Range("J138").Value = 4036.5293
$ws.Range("L138").Value = 12109.5879
$ws.Range("N138").Value = -22389.5879

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 900.5454999999999
$ws.Range("I2").Value = 906.1667
$ws.Range("J2").Value = 893.8
$ws.Range("K2").Value = 906.1667
$ws.Range("L2").Value = 893.8
$ws.Range("M2").Value = -793.1667
$ws.Range("N2").Value = -1119.8
$ws.Range("H61").Value = 1901.2
$ws.Range("I61").Value = 1834.8334
$ws.Range("J61").Value = 2498.5
$ws.Range("K61").Value = 1834.8334
$ws.Range("L61").Value = 2498.5
$ws.Range("M61").Value = -1622.8334
$ws.Range("N61").Value = -2922.5
$ws.Range("H102").Value = 22996.143
$ws.Range("I102").Value = 21599.2
$ws.Range("K102").Value = 21599.2
$ws.Range("M102").Value = -19977.2
$ws.Range("H116").Value = 900.5454999999999
$ws.Range("I116").Value = 906.1667
$ws.Range("J116").Value = 893.8
$ws.Range("K116").Value = 906.1667
$ws.Range("L116").Value = 893.8
$ws.Range("M116").Value = 1387.8333
$ws.Range("N116").Value = -5481.8
$ws.Range("H136").Value = 1901.2
$ws.Range("I136").Value = 1834.8334
$ws.Range("J136").Value = 2498.5
$ws.Range("K136").Value = 5504.5002
$ws.Range("L136").Value = 7495.5
$ws.Range("M136").Value = -2954.5002
$ws.Range("N136").Value = -12595.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 900.5454999999999
$ws.Range("I3").Value = 906.1667
$ws.Range("J3").Value = 893.8
$ws.Range("K3").Value = 906.1667
$ws.Range("L3").Value = 893.8
$ws.Range("M3").Value = -792.1667
$ws.Range("N3").Value = -1121.8
$ws.Range("H94").Value = 7099.778
$ws.Range("I94").Value = 8780.143
$ws.Range("K94").Value = 8780.143
$ws.Range("M94").Value = -8329.143

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 7002
$ws.Range("I2").Value = 4004
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 4004
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = -3891
$ws.Range("N2").Value = -10226
$ws.Range("H31").Value = 5279.1816
$ws.Range("I31").Value = 2966
$ws.Range("K31").Value = 2966
$ws.Range("M31").Value = -2671
$ws.Range("H34").Value = 5279.1816
$ws.Range("I34").Value = 2966
$ws.Range("K34").Value = 2966
$ws.Range("M34").Value = -2764
$ws.Range("H99").Value = 1380.091
$ws.Range("J99").Value = 846
$ws.Range("L99").Value = 846
$ws.Range("N99").Value = -3842
$ws.Range("H105").Value = 1499
$ws.Range("I105").Value = 1499
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1499
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 248
$ws.Range("N105").Value = $null
$ws.Range("H107").Value = 774.5454999999999
$ws.Range("I107").Value = 408.44446
$ws.Range("K107").Value = 408.44446
$ws.Range("M107").Value = 1511.55554
$ws.Range("H122").Value = 5000
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900
$ws.Range("H126").Value = 1380.091
$ws.Range("J126").Value = 846
$ws.Range("L126").Value = 2538
$ws.Range("N126").Value = -7478

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 290.54544
$ws.Range("J12").Value = 241.2
$ws.Range("L12").Value = 723.5999999999999
$ws.Range("N12").Value = -1069.6
$ws.Range("H26").Value = 749.8182
$ws.Range("I26").Value = 599.5
$ws.Range("J26").Value = 835.7143
$ws.Range("K26").Value = 1798.5
$ws.Range("L26").Value = 2507.1429
$ws.Range("M26").Value = -1510.5
$ws.Range("N26").Value = -3083.1429
$ws.Range("H132").Value = 4696.5
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 5107.222
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 45964.998
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -51024.998

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1092.5454
$ws.Range("I113").Value = 1113.1111
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1113.1111
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1056.8889
$ws.Range("N113").Value = -5340
$ws.Range("H126").Value = 1954.8
$ws.Range("I126").Value = 1387.25
$ws.Range("K126").Value = 4161.75
$ws.Range("M126").Value = -1691.75
$ws.Range("H132").Value = 3879.2
$ws.Range("I132").Value = 3553.077
$ws.Range("K132").Value = 10659.231
$ws.Range("M132").Value = -8129.231

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3716
$ws.Range("I46").Value = 3696.3333
$ws.Range("J46").Value = 3775
$ws.Range("K46").Value = 3696.3333
$ws.Range("L46").Value = 3775
$ws.Range("M46").Value = -3508.3333
$ws.Range("N46").Value = -4151
$ws.Range("H61").Value = 4701.636
$ws.Range("I61").Value = 4749.5713
$ws.Range("J61").Value = 4617.75
$ws.Range("K61").Value = 4749.5713
$ws.Range("L61").Value = 4617.75
$ws.Range("M61").Value = -4547.5713
$ws.Range("N61").Value = -5021.75
$ws.Range("H82").Value = 566.25
$ws.Range("I82").Value = 605
$ws.Range("J82").Value = 295
$ws.Range("K82").Value = 605
$ws.Range("L82").Value = 295
$ws.Range("M82").Value = -244
$ws.Range("N82").Value = -1017
$ws.Range("H85").Value = 566.25
$ws.Range("I85").Value = 605
$ws.Range("J85").Value = 295
$ws.Range("K85").Value = 605
$ws.Range("L85").Value = 295
$ws.Range("M85").Value = 643
$ws.Range("N85").Value = -2791
$ws.Range("H113").Value = 4701.636
$ws.Range("I113").Value = 4749.5713
$ws.Range("J113").Value = 4617.75
$ws.Range("K113").Value = 4749.5713
$ws.Range("L113").Value = 4617.75
$ws.Range("M113").Value = -2579.5713
$ws.Range("N113").Value = -8957.75

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26498
$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -82488
$ws.Range("H132").Value = 2836.258
$ws.Range("I132").Value = 2478.2693
$ws.Range("K132").Value = 7434.8079
$ws.Range("M132").Value = -4904.8079
$ws.Range("H136").Value = 7039.7915
$ws.Range("I136").Value = 7039.7915
$ws.Range("K136").Value = 21119.3745
$ws.Range("M136").Value = -18569.3745
